# Add "Week 3" time-log entries to each team member's sheet and roll the
# totals up into a new "Week 3" block on the Overview sheet.

$wb = $excel.ActiveWorkbook

# --- Per-person sheets: append a "Week 3" row with that week's hours -------
$weekly = @{
    "Nedas_J"      = 0.54166666666666663
    "Adomas_J"     = 0.17361111111111113
    "Aistė_G"      = 0.25347222222222221
    "Gabrielius_D" = 0.18402777777777779
}

foreach ($name in @("Nedas_J", "Adomas_J", "Aistė_G", "Gabrielius_D")) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A3").Value = "Week 3"
    $ws.Range("B3").NumberFormat = $ws.Range("B2").NumberFormat
    $ws.Range("B3").Value = $weekly[$name]
}

# --- Overview sheet: new "Week 3 (working in a group)" / "(total working hours)" rows
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A5").Value = "Week 3 (working in a group)"
$ov.Range("B3").Copy() | Out-Null
$ov.Range("B5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ov.Range("B5").Value = 0

$ov.Range("A6").Value = "Week 3 (total working hours)"
$ov.Range("B4").Copy() | Out-Null
$ov.Range("B6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ov.Range("B6").Formula = "=SUM(Nedas_J!B3, Adomas_J!B3, Aistė_G!B3, Gabrielius_D!B3) + 4 * B5"

# --- Restore the on-screen selections Excel leaves behind after such an edit
$wb.Worksheets.Item("Adomas_J").Range("B4").Select() | Out-Null
$wb.Worksheets.Item("Aistė_G").Range("B4").Select() | Out-Null
$wb.Worksheets.Item("Gabrielius_D").Range("F30").Select() | Out-Null

$ov.Select() | Out-Null
$ov.Range("B6").Select() | Out-Null
$excel.ActiveWindow.Zoom = 189

$wb.Worksheets.Item("Nedas_J").Range("B4").Select() | Out-Null
